$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status by State")

# Combine column C ("New") and column D ("Present") into column C, then drop column D.
for ($r = 2; $r -le 39; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    if ($cVal -eq $null) { $cVal = 0 }
    if ($dVal -eq $null) { $dVal = 0 }
    $ws.Cells.Item($r, 3).Value = $cVal + $dVal
}

# Header: C1 becomes "Present" (what D1 used to say)
$ws.Cells.Item(1, 3).Value = "Present"

# Remove column D entirely
$ws.Columns.Item(4).Delete()
